$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.988074333333333
$ws.Range("H2").Value = 5.964223
$ws.Range("I2").Value = 0.01657769708907969
$ws.Range("J2").Value = 0.01657769708907968
$ws.Range("M2").Value = 1.174933333333333
$ws.Range("N2").Value = 3.5248
$ws.Range("O2").Value = 0.01171850713626266
$ws.Range("P2").Value = 0.01171850713626266
$ws.Range("Q2").Value = 2.335854803377778
$ws.Range("R2").Value = 21.0226932304
$ws.Range("S2").Value = 0.000194265861641181
$ws.Range("T2").Value = 0.0001942658616411809
$ws.Range("G3").Value = 1.988074333333333
$ws.Range("H3").Value = 5.964223
$ws.Range("I3").Value = 0.01657769708907969
$ws.Range("J3").Value = 0.01657769708907968
$ws.Range("O3").Value = 0.2743256641287217
$ws.Range("P3").Value = 0.2743256641287218
$ws.Range("Q3").Value = 54.68144643288044
$ws.Range("R3").Value = 492.133017895924
$ws.Range("S3").Value = 0.004547687763686562
$ws.Range("T3").Value = 0.004547687763686562
$ws.Range("G4").Value = 1.988074333333333
$ws.Range("H4").Value = 5.964223
$ws.Range("I4").Value = 0.01657769708907969
$ws.Range("J4").Value = 0.01657769708907968
$ws.Range("M4").Value = 39.361408
$ws.Range("N4").Value = 118.084224
$ws.Range("O4").Value = 0.3925813724534833
$ws.Range("P4").Value = 0.3925813724534833
$ws.Range("Q4").Value = 78.25340496866134
$ws.Range("R4").Value = 704.280644717952
$ws.Range("S4").Value = 0.006508095075349019
$ws.Range("T4").Value = 0.006508095075349018
$ws.Range("G5").Value = 1.988074333333333
$ws.Range("H5").Value = 5.964223
$ws.Range("I5").Value = 0.01657769708907969
$ws.Range("J5").Value = 0.01657769708907968
$ws.Range("M5").Value = 32.221985
$ws.Range("N5").Value = 96.665955
$ws.Range("O5").Value = 0.3213744562815322
$ws.Range("P5").Value = 0.3213744562815322
$ws.Range("Q5").Value = 64.05970134755165
$ws.Range("R5").Value = 576.5373121279649
$ws.Range("S5").Value = 0.005327648388402923
$ws.Range("T5").Value = 0.005327648388402923
$ws.Range("I6").Value = 0.7746030815641455
$ws.Range("J6").Value = 0.7746030815641454
$ws.Range("M6").Value = 1.174933333333333
$ws.Range("N6").Value = 3.5248
$ws.Range("O6").Value = 0.01171850713626266
$ws.Range("P6").Value = 0.01171850713626266
$ws.Range("Q6").Value = 109.1442508003556
$ws.Range("R6").Value = 982.2982572032
$ws.Range("S6").Value = 0.009077191739080484
$ws.Range("T6").Value = 0.009077191739080481
$ws.Range("I7").Value = 0.7746030815641455
$ws.Range("J7").Value = 0.7746030815641454
$ws.Range("O7").Value = 0.2743256641287217
$ws.Range("P7").Value = 0.2743256641287218
$ws.Range("S7").Value = 0.2124935047862386
$ws.Range("T7").Value = 0.2124935047862386
$ws.Range("I8").Value = 0.7746030815641455
$ws.Range("J8").Value = 0.7746030815641454
$ws.Range("M8").Value = 39.361408
$ws.Range("N8").Value = 118.084224
$ws.Range("O8").Value = 0.3925813724534833
$ws.Range("P8").Value = 0.3925813724534833
$ws.Range("Q8").Value = 3656.438424824491
$ws.Range("R8").Value = 32907.94582342042
$ws.Range("S8").Value = 0.3040947408671497
$ws.Range("T8").Value = 0.3040947408671497
$ws.Range("I9").Value = 0.7746030815641455
$ws.Range("J9").Value = 0.7746030815641454
$ws.Range("M9").Value = 32.221985
$ws.Range("N9").Value = 96.665955
$ws.Range("O9").Value = 0.3213744562815322
$ws.Range("P9").Value = 0.3213744562815322
$ws.Range("Q9").Value = 2993.228902739413
$ws.Range("R9").Value = 26939.06012465472
$ws.Range("S9").Value = 0.2489376441716766
$ws.Range("T9").Value = 0.2489376441716766
$ws.Range("G10").Value = 23.741365
$ws.Range("H10").Value = 71.224095
$ws.Range("I10").Value = 0.1979690350870239
$ws.Range("J10").Value = 0.1979690350870239
$ws.Range("M10").Value = 1.174933333333333
$ws.Range("N10").Value = 3.5248
$ws.Range("O10").Value = 0.01171850713626266
$ws.Range("P10").Value = 0.01171850713626266
$ws.Range("Q10").Value = 27.89452111733334
$ws.Range("R10").Value = 251.050690056
$ws.Range("S10").Value = 0.002319901550426322
$ws.Range("T10").Value = 0.002319901550426322
$ws.Range("G11").Value = 23.741365
$ws.Range("H11").Value = 71.224095
$ws.Range("I11").Value = 0.1979690350870239
$ws.Range("J11").Value = 0.1979690350870239
$ws.Range("O11").Value = 0.2743256641287217
$ws.Range("P11").Value = 0.2743256641287218
$ws.Range("Q11").Value = 652.9998183288734
$ws.Range("R11").Value = 5876.998364959861
$ws.Range("S11").Value = 0.05430798702717005
$ws.Range("T11").Value = 0.05430798702717007
$ws.Range("G12").Value = 23.741365
$ws.Range("H12").Value = 71.224095
$ws.Range("I12").Value = 0.1979690350870239
$ws.Range("J12").Value = 0.1979690350870239
$ws.Range("M12").Value = 39.361408
$ws.Range("N12").Value = 118.084224
$ws.Range("O12").Value = 0.3925813724534833
$ws.Range("P12").Value = 0.3925813724534833
$ws.Range("Q12").Value = 934.4935542419202
$ws.Range("R12").Value = 8410.441988177281
$ws.Range("S12").Value = 0.07771895549775565
$ws.Range("T12").Value = 0.07771895549775565
$ws.Range("G13").Value = 23.741365
$ws.Range("H13").Value = 71.224095
$ws.Range("I13").Value = 0.1979690350870239
$ws.Range("J13").Value = 0.1979690350870239
$ws.Range("M13").Value = 32.221985
$ws.Range("N13").Value = 96.665955
$ws.Range("O13").Value = 0.3213744562815322
$ws.Range("P13").Value = 0.3213744562815322
$ws.Range("Q13").Value = 764.993906909525
$ws.Range("R13").Value = 6884.945162185725
$ws.Range("S13").Value = 0.06362219101167188
$ws.Range("T13").Value = 0.0636221910116719
$ws.Range("G14").Value = 1.301204666666667
$ws.Range("H14").Value = 3.903614
$ws.Range("I14").Value = 0.01085018625975097
$ws.Range("J14").Value = 0.01085018625975097
$ws.Range("M14").Value = 1.174933333333333
$ws.Range("N14").Value = 3.5248
$ws.Range("O14").Value = 0.01171850713626266
$ws.Range("P14").Value = 0.01171850713626266
$ws.Range("Q14").Value = 1.528828736355556
$ws.Range("R14").Value = 13.7594586272
$ws.Range("S14").Value = 0.0001271479851146708
$ws.Range("T14").Value = 0.0001271479851146708
$ws.Range("G15").Value = 1.301204666666667
$ws.Range("H15").Value = 3.903614
$ws.Range("I15").Value = 0.01085018625975097
$ws.Range("J15").Value = 0.01085018625975097
$ws.Range("O15").Value = 0.2743256641287217
$ws.Range("P15").Value = 0.2743256641287218
$ws.Range("Q15").Value = 35.78928216393689
$ws.Range("R15").Value = 322.103539475432
$ws.Range("S15").Value = 0.002976484551626517
$ws.Range("T15").Value = 0.002976484551626517
$ws.Range("G16").Value = 1.301204666666667
$ws.Range("H16").Value = 3.903614
$ws.Range("I16").Value = 0.01085018625975097
$ws.Range("J16").Value = 0.01085018625975097
$ws.Range("M16").Value = 39.361408
$ws.Range("N16").Value = 118.084224
$ws.Range("O16").Value = 0.3925813724534833
$ws.Range("P16").Value = 0.3925813724534833
$ws.Range("Q16").Value = 51.21724777617068
$ws.Range("R16").Value = 460.9552299855361
$ws.Range("S16").Value = 0.004259581013228964
$ws.Range("T16").Value = 0.004259581013228963
$ws.Range("G17").Value = 1.301204666666667
$ws.Range("H17").Value = 3.903614
$ws.Range("I17").Value = 0.01085018625975097
$ws.Range("J17").Value = 0.01085018625975097
$ws.Range("M17").Value = 32.221985
$ws.Range("N17").Value = 96.665955
$ws.Range("O17").Value = 0.3213744562815322
$ws.Range("P17").Value = 0.3213744562815322
$ws.Range("Q17").Value = 41.92739725126333
$ws.Range("R17").Value = 377.34657526137
$ws.Range("S17").Value = 0.00348697270978082
$ws.Range("T17").Value = 0.00348697270978082
